# Generate Report for Handback
# Updates the "Status" text for the b7d19e9e handoff row (it failed the
# handback transform), records the detailed error message in the
# "Error Detail" column of the zh-cn / de-de sheets, and widens that
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$zhError = "Handback file name: x5nnnvfs.yue is different with handoff file name: b7d19e9e-8016-4caf-8543-fb97cdbc792d.337bdecf6b408406ccf4f644deeb252ac3a9c38c.zh-cn."
$deError = "Handback file name: x5nnnvfs.yue is different with handoff file name: b7d19e9e-8016-4caf-8543-fb97cdbc792d.337bdecf6b408406ccf4f644deeb252ac3a9c38c.de-de."

# --- Overview sheet: Status column (E/F) for the b7d19e9e row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C) row 3, Error Detail column (P) row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status column (C) row 3, Error Detail column (P) row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17
